$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.264697432518005
$ws.Range("B1").Value = 1.365736722946167
$ws.Range("C1").Value = 1.541692137718201
$ws.Range("D1").Value = 2.402352094650269
$ws.Range("E1").Value = 15
